# Edit: bump "Chapter 01" heading to "Chapter 02" (keeping the bookmark end
# anchored right after "Chapter 0") and insert a new notes table (with the
# "Page 17" Angular/TypeScript note) after the chapter heading, before the
# existing empty "Titulo 3" paragraph.

$d = $word.ActiveDocument

# --- Step 1: "Chapter 01" -> "Chapter 0" / bookmarkEnd / "2" -------------
$chapterBookmark = $d.Bookmarks("_Toc147398633")
$lastDigit = $d.Range($chapterBookmark.End - 1, $chapterBookmark.End)
$lastDigit.Delete()
$insertPoint = $d.Range($chapterBookmark.End, $chapterBookmark.End)
$insertPoint.InsertAfter("2")

# --- Step 2: insert blank paragraph + notes table + blank paragraph ------
# right after the chapter heading paragraph (before the following, empty
# "Titulo 3" paragraph).
$chapterPara = $chapterBookmark.Range.Paragraphs(1)
$insertAt = $chapterPara.Range.End - 1
$tblRange = $d.Range($insertAt, $insertAt)
$tableXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:tbl><w:tblPr><w:tblStyle w:val="TabelacomGrelha"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="9016"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Page 17:</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">Angular applications are written in TypeScript, which is a superset of JavaScript. I introduce TypeScript in Chapters 3 and 4, but its main advantage is that it supports </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>static data types</w:t></w:r><w:r><w:t xml:space="preserve">, which makes JavaScript development more familiar to C# and Java developers. (JavaScript has a </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>prototype-based type system</w:t></w:r><w:r><w:t xml:space="preserve"> that many developers find confusing.)</w:t></w:r></w:p><w:p/><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/></w:tcPr><w:p/></w:tc></w:tr></w:tbl><w:p/>'
[void]$tblRange.InsertXML($tableXml)
